# Scheduled runner update: refresh market-board derived price/profit columns
# (currentAveragePrice, currentAveragePriceNQ, currentAveragePriceHQ,
#  LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ -> columns H:N)
# for the affected leve rows across the ALC, ARM, CRP, CUL, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3936.5386
$ws.Range("I76").Value = 3898.889
$ws.Range("J76").Value = 4021.25
$ws.Range("K76").Value = 3898.889
$ws.Range("L76").Value = 4021.25
$ws.Range("M76").Value = -3583.889
$ws.Range("N76").Value = -4651.25

$ws.Range("H79").Value = 3936.5386
$ws.Range("I79").Value = 3898.889
$ws.Range("J79").Value = 4021.25
$ws.Range("K79").Value = 3898.889
$ws.Range("L79").Value = 4021.25
$ws.Range("M79").Value = -2806.889
$ws.Range("N79").Value = -6205.25

$ws.Range("H86").Value = 1090.2858
$ws.Range("I86").Value = 1092.6364
$ws.Range("J86").Value = 1081.6666
$ws.Range("K86").Value = 1092.6364
$ws.Range("L86").Value = 1081.6666
$ws.Range("M86").Value = 30.36359999999991
$ws.Range("N86").Value = -3327.6666

$ws.Range("H89").Value = 1090.2858
$ws.Range("I89").Value = 1092.6364
$ws.Range("J89").Value = 1081.6666
$ws.Range("K89").Value = 5463.182000000001
$ws.Range("L89").Value = 5408.333000000001
$ws.Range("M89").Value = 152.8179999999993
$ws.Range("N89").Value = -16640.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 819.39624
$ws.Range("I74").Value = 686.7143
$ws.Range("J74").Value = 1326
$ws.Range("K74").Value = 686.7143
$ws.Range("L74").Value = 1326
$ws.Range("M74").Value = 187.2857
$ws.Range("N74").Value = -3074

$ws.Range("H77").Value = 819.39624
$ws.Range("I77").Value = 686.7143
$ws.Range("J77").Value = 1326
$ws.Range("K77").Value = 3433.5715
$ws.Range("L77").Value = 6630
$ws.Range("M77").Value = 934.4285
$ws.Range("N77").Value = -15366

$ws.Range("H132").Value = 2393.8108
$ws.Range("I132").Value = 1771.5
$ws.Range("J132").Value = 2772.6086
$ws.Range("K132").Value = 5314.5
$ws.Range("L132").Value = 8317.825800000001
$ws.Range("M132").Value = -2784.5
$ws.Range("N132").Value = -13377.8258

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2082
$ws.Range("I31").Value = 1559.6666
$ws.Range("J31").Value = 8350
$ws.Range("K31").Value = 1559.6666
$ws.Range("L31").Value = 8350
$ws.Range("M31").Value = -1264.6666
$ws.Range("N31").Value = -8940

$ws.Range("H34").Value = 2082
$ws.Range("I34").Value = 1559.6666
$ws.Range("J34").Value = 8350
$ws.Range("K34").Value = 1559.6666
$ws.Range("L34").Value = 8350
$ws.Range("M34").Value = -1357.6666
$ws.Range("N34").Value = -8754

$ws.Range("H58").Value = 883120.7
$ws.Range("I58").Value = 1235625.9
$ws.Range("J58").Value = 1857.75
$ws.Range("K58").Value = 1235625.9
$ws.Range("L58").Value = 1857.75
$ws.Range("M58").Value = -1235422.9
$ws.Range("N58").Value = -2263.75

$ws.Range("H99").Value = 3612.9285
$ws.Range("I99").Value = 4074.111
$ws.Range("J99").Value = 2782.8
$ws.Range("K99").Value = 4074.111
$ws.Range("L99").Value = 2782.8
$ws.Range("M99").Value = -2576.111
$ws.Range("N99").Value = -5778.8

$ws.Range("H126").Value = 3612.9285
$ws.Range("I126").Value = 4074.111
$ws.Range("J126").Value = 2782.8
$ws.Range("K126").Value = 12222.333
$ws.Range("L126").Value = 8348.400000000001
$ws.Range("M126").Value = -9752.332999999999
$ws.Range("N126").Value = -13288.4

$ws.Range("H132").Value = 196855.75
$ws.Range("I132").Value = 246439.97
$ws.Range("J132").Value = 2060.6428
$ws.Range("K132").Value = 739319.91
$ws.Range("L132").Value = 6181.928400000001
$ws.Range("M132").Value = -736789.91
$ws.Range("N132").Value = -11241.9284

$ws.Range("H134").Value = 1129.5903
$ws.Range("I134").Value = 1013.3043
$ws.Range("J134").Value = 1702.7142
$ws.Range("K134").Value = 3039.9129
$ws.Range("L134").Value = 5108.142599999999
$ws.Range("M134").Value = -504.9129000000003
$ws.Range("N134").Value = -10178.1426

$ws.Range("H136").Value = 883120.7
$ws.Range("I136").Value = 1235625.9
$ws.Range("J136").Value = 1857.75
$ws.Range("K136").Value = 3706877.7
$ws.Range("L136").Value = 5573.25
$ws.Range("M136").Value = -3704327.7
$ws.Range("N136").Value = -10673.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2129.8667
$ws.Range("I5").Value = 2363.4546
$ws.Range("K5").Value = 7090.3638
$ws.Range("M5").Value = -6978.3638

$ws.Range("H75").Value = 7691.5
$ws.Range("I75").Value = 1633
$ws.Range("J75").Value = 13750
$ws.Range("K75").Value = 4899
$ws.Range("L75").Value = 41250
$ws.Range("M75").Value = -3901
$ws.Range("N75").Value = -43246

$ws.Range("H78").Value = 7691.5
$ws.Range("I78").Value = 1633
$ws.Range("J78").Value = 13750
$ws.Range("K78").Value = 14697
$ws.Range("L78").Value = 123750
$ws.Range("M78").Value = -9705
$ws.Range("N78").Value = -133734

$ws.Range("H129").Value = 1667918.4
$ws.Range("I129").Value = 667.3333
$ws.Range("J129").Value = 2779419.2
$ws.Range("K129").Value = 2001.9999
$ws.Range("L129").Value = 8338257.600000001
$ws.Range("M129").Value = 2998.0001
$ws.Range("N129").Value = -8348257.600000001

$ws.Range("H131").Value = 23813880
$ws.Range("I131").Value = 22186
$ws.Range("J131").Value = 27028974
$ws.Range("K131").Value = 66558
$ws.Range("L131").Value = 81086922
$ws.Range("M131").Value = -61518
$ws.Range("N131").Value = -81097002

$ws.Range("H135").Value = 2129.8667
$ws.Range("I135").Value = 2363.4546
$ws.Range("K135").Value = 21271.0914
$ws.Range("M135").Value = -18736.0914

$ws.Range("H137").Value = 20836966
$ws.Range("I137").Value = 3106
$ws.Range("J137").Value = 30306902
$ws.Range("K137").Value = 9318
$ws.Range("L137").Value = 90920706
$ws.Range("M137").Value = -4218
$ws.Range("N137").Value = -90930906

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1016.1539
$ws.Range("I61").Value = 1064.5454
$ws.Range("J61").Value = 750
$ws.Range("K61").Value = 1064.5454
$ws.Range("L61").Value = 750
$ws.Range("M61").Value = -862.5454
$ws.Range("N61").Value = -1154

$ws.Range("H113").Value = 1016.1539
$ws.Range("I113").Value = 1064.5454
$ws.Range("J113").Value = 750
$ws.Range("K113").Value = 1064.5454
$ws.Range("L113").Value = 750
$ws.Range("M113").Value = 1105.4546
$ws.Range("N113").Value = -5090

$ws.Range("H132").Value = 2707.75
$ws.Range("I132").Value = 2217.0715
$ws.Range("J132").Value = 4997.5835
$ws.Range("K132").Value = 6651.2145
$ws.Range("L132").Value = 14992.7505
$ws.Range("M132").Value = -4121.2145
$ws.Range("N132").Value = -20052.7505

$ws.Range("H136").Value = 2688.7793
$ws.Range("I136").Value = 2653.1775
$ws.Range("J136").Value = 2835.9333
$ws.Range("K136").Value = 7959.532499999999
$ws.Range("L136").Value = 8507.7999
$ws.Range("M136").Value = -5409.532499999999
$ws.Range("N136").Value = -13607.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 917.0328
$ws.Range("I132").Value = 655.7805
$ws.Range("J132").Value = 1452.6
$ws.Range("K132").Value = 1967.3415
$ws.Range("L132").Value = 4357.799999999999
$ws.Range("M132").Value = 562.6585
$ws.Range("N132").Value = -9417.799999999999

$ws.Range("H136").Value = 1753.7826
$ws.Range("I136").Value = 1386.1578
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 4158.4734
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -1608.4734
$ws.Range("N136").Value = -15600
